$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 17599.88
$ws.Range("I43").Value = 17166.611
$ws.Range("J43").Value = 18714
$ws.Range("K43").Value = 17166.611
$ws.Range("L43").Value = 18714
$ws.Range("M43").Value = -17097.611
$ws.Range("N43").Value = -18852
$ws.Range("H111").Value = 759.36365
$ws.Range("I111").Value = 604.8570999999999
$ws.Range("J111").Value = 1029.75
$ws.Range("K111").Value = 1814.5713
$ws.Range("L111").Value = 3089.25
$ws.Range("M111").Value = 1252.4287
$ws.Range("N111").Value = -9223.25
$ws.Range("H132").Value = 1512.9796
$ws.Range("I132").Value = 1501.3864
$ws.Range("K132").Value = 4504.1592
$ws.Range("M132").Value = -1974.1592
$ws.Range("H137").Value = 3632.0967
$ws.Range("J137").Value = 4045.4736
$ws.Range("L137").Value = 12136.4208
$ws.Range("N137").Value = -17236.4208
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2001.5714
$ws.Range("I32").Value = 1836.0597
$ws.Range("J32").Value = 5698
$ws.Range("K32").Value = 1836.0597
$ws.Range("L32").Value = 5698
$ws.Range("M32").Value = -1549.0597
$ws.Range("N32").Value = -6272
$ws.Range("H63").Value = 1947.5
$ws.Range("I63").Value = 1947.5
$ws.Range("K63").Value = 1947.5
$ws.Range("M63").Value = -1261.5
$ws.Range("H66").Value = 1947.5
$ws.Range("I66").Value = 1947.5
$ws.Range("K66").Value = 9737.5
$ws.Range("M66").Value = -6305.5
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H110").Value = 7949.75
$ws.Range("I110").Value = 7899.5
$ws.Range("K110").Value = 7899.5
$ws.Range("M110").Value = -5854.5
$ws.Range("H119").Value = 68127
$ws.Range("J119").Value = 68127
$ws.Range("L119").Value = 68127
$ws.Range("N119").Value = -77803
$ws.Range("H132").Value = 9657.879999999999
$ws.Range("J132").Value = 15378.182
$ws.Range("L132").Value = 46134.546
$ws.Range("N132").Value = -51194.546
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1487.4445
$ws.Range("I134").Value = 967.08
$ws.Range("K134").Value = 2901.24
$ws.Range("M134").Value = -366.2400000000002
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10677.875
$ws.Range("I62").Value = 4569.5
$ws.Range("J62").Value = 29003
$ws.Range("K62").Value = 4569.5
$ws.Range("L62").Value = 29003
$ws.Range("M62").Value = -3945.5
$ws.Range("N62").Value = -30251
$ws.Range("H65").Value = 10677.875
$ws.Range("I65").Value = 4569.5
$ws.Range("J65").Value = 29003
$ws.Range("K65").Value = 22847.5
$ws.Range("L65").Value = 145015
$ws.Range("M65").Value = -19727.5
$ws.Range("N65").Value = -151255
$ws.Range("H105").Value = 5273.5
$ws.Range("I105").Value = 5098.143
$ws.Range("K105").Value = 5098.143
$ws.Range("M105").Value = -3351.143
$ws.Range("H122").Value = 51645.273
$ws.Range("I122").Value = 65494
$ws.Range("J122").Value = 14715.333
$ws.Range("K122").Value = 196482
$ws.Range("L122").Value = 44145.999
$ws.Range("M122").Value = -194032
$ws.Range("N122").Value = -49045.999
$ws.Range("H125").Value = 133108.67
$ws.Range("J125").Value = 133108.67
$ws.Range("L125").Value = 133108.67
$ws.Range("N125").Value = -138028.67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4323408.5
$ws.Range("I34").Value = 6482679
$ws.Range("J34").Value = 4866.5
$ws.Range("K34").Value = 19448037
$ws.Range("L34").Value = 14599.5
$ws.Range("M34").Value = -19447953
$ws.Range("N34").Value = -14767.5
$ws.Range("H39").Value = 2485.1428
$ws.Range("J39").Value = 2732
$ws.Range("L39").Value = 8196
$ws.Range("N39").Value = -8784
$ws.Range("H55").Value = 1879.1111
$ws.Range("I55").Value = 1382.2
$ws.Range("J55").Value = 2500.25
$ws.Range("K55").Value = 4146.6
$ws.Range("L55").Value = 7500.75
$ws.Range("M55").Value = -3969.6
$ws.Range("N55").Value = -7854.75
$ws.Range("H92").Value = 998.1539
$ws.Range("J92").Value = 2368.25
$ws.Range("L92").Value = 7104.75
$ws.Range("N92").Value = -9600.75
$ws.Range("H122").Value = 1167461.1
$ws.Range("J122").Value = 2386.2354
$ws.Range("L122").Value = 21476.1186
$ws.Range("N122").Value = -26376.1186
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 835900.3
$ws.Range("I3").Value = 1667283.4
$ws.Range("J3").Value = 4517.3335
$ws.Range("K3").Value = 1667283.4
$ws.Range("L3").Value = 4517.3335
$ws.Range("M3").Value = -1667167.4
$ws.Range("N3").Value = -4749.3335
$ws.Range("H24").Value = 53906
$ws.Range("I24").Value = 53906
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 53906
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -53733
$ws.Range("N24").ClearContents()
$ws.Range("H122").Value = 4983.552
$ws.Range("I122").Value = 5161.7827
$ws.Range("J122").Value = 4300.3335
$ws.Range("K122").Value = 15485.3481
$ws.Range("L122").Value = 12901.0005
$ws.Range("M122").Value = -13035.3481
$ws.Range("N122").Value = -17801.0005
$ws.Range("H132").Value = 489570.34
$ws.Range("I132").Value = 627808.8
$ws.Range("J132").Value = 5735.6665
$ws.Range("K132").Value = 1883426.4
$ws.Range("L132").Value = 17206.9995
$ws.Range("M132").Value = -1880896.4
$ws.Range("N132").Value = -22266.9995
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6393.7427
$ws.Range("I7").Value = 3906.5557
$ws.Range("J7").Value = 14788
$ws.Range("K7").Value = 3906.5557
$ws.Range("L7").Value = 14788
$ws.Range("M7").Value = -3794.5557
$ws.Range("N7").Value = -15012
$ws.Range("H40").Value = 8525.700000000001
$ws.Range("I40").Value = 9139.666999999999
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 9139.666999999999
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -9003.666999999999
$ws.Range("N40").Value = -3272
$ws.Range("H50").Value = 30750
$ws.Range("J50").Value = 30750
$ws.Range("L50").Value = 30750
$ws.Range("N50").Value = -32024
$ws.Range("H61").Value = 3128.85
$ws.Range("I61").Value = 1462.5
$ws.Range("J61").Value = 7017
$ws.Range("K61").Value = 1462.5
$ws.Range("L61").Value = 7017
$ws.Range("M61").Value = -1260.5
$ws.Range("N61").Value = -7421
$ws.Range("H113").Value = 3128.85
$ws.Range("I113").Value = 1462.5
$ws.Range("J113").Value = 7017
$ws.Range("K113").Value = 1462.5
$ws.Range("L113").Value = 7017
$ws.Range("M113").Value = 707.5
$ws.Range("N113").Value = -11357
$ws.Range("H122").Value = 6700.077
$ws.Range("I122").Value = 4009.6
$ws.Range("K122").Value = 12028.8
$ws.Range("M122").Value = -9578.799999999999
$ws.Range("H126").Value = 6393.7427
$ws.Range("I126").Value = 3906.5557
$ws.Range("J126").Value = 14788
$ws.Range("K126").Value = 11719.6671
$ws.Range("L126").Value = 44364
$ws.Range("M126").Value = -9249.667099999999
$ws.Range("N126").Value = -49304
$ws.Range("H140").Value = 90429
$ws.Range("J140").Value = 90429
$ws.Range("L140").Value = 90429
$ws.Range("N140").Value = -100789
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 16681667
$ws.Range("J5").Value = 25010000
$ws.Range("L5").Value = 25010000
$ws.Range("N5").Value = -25010224
$ws.Range("H132").Value = 6557.9707
$ws.Range("I132").Value = 2658.682
$ws.Range("K132").Value = 7976.045999999999
$ws.Range("M132").Value = -5446.045999999999
